# [zvnl 6400] scenario Ehv-Tbu
# Adds a new "6400 Ehv-Tb" scenario column (C) to the "QD" sheet, splits /
# renames some of the existing "4400 Ehv-Ht" scenario notes in column B,
# and refreshes the view-selection state.
#
# NOTE: the order of the `.Value =` assignments below is deliberate — it
# reproduces the exact shared-string insertion order of the original edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QD")

# 1) New header cell for the new scenario column.
$ws.Range("C1").Value = "6400 Ehv-Tb"

# 2) Timers note for the new scenario.
$ws.Range("C13").Value = "- Ehv"

# 3) AI note for the new scenario.
$ws.Range("C4").Value = "Ehv:`n- geen"

# 4) B5's old combined note is split: the tail ("- 3900 Btl / - G Btl-Vg /
#    - 3500 Vga / - 4400 Vga") becomes its own "Btl-Ht:" note in B6 ...
$ws.Range("B6").Value = "Btl-Ht:`n- 3900 Btl`n- G Btl-Vg`n- 3500 Vga`n- 4400 Vga"

# 5) ... and B5 keeps the head of the note, renamed "Ehv-Btl:".
$ws.Range("B5").Value = "Ehv-Btl:`n- 6400 Ehs`n- G Ehs`n- 3500 At`n- 1100 At`n- G 2x Bet`n- 4400 Bet`n- 800 Beto`n- 1900 Beto`n- G Beto-Lpe`n- G Lpe`n- 6400 Btl"

# 6) New column C counterpart notes for the new scenario.
$ws.Range("C5").Value = "Ehb-Btl:`n- 800 Ehs`n- 4400 Ehs`n- G 2x At`n- 3900 Bet`n- 6400 Bet`n- G 2x Beto`n- 3500 Lpe`n- 1100 Lpe`n- 4400 Btl"
$ws.Range("C6").Value = "Btl-Tbu:`n- G Btl-Otw`n- 6400 Otw`n- G Otw-Tba`n- 1100 Tb`n- G Tbu"
$ws.Range("C12").Value = "- G voor sp`n- 3900 achter sp"
$ws.Range("C11").Value = "- At"

# Old B6 note ("Ht: ...") shifts down to B7 (text unchanged, reuses the
# existing shared string).
$ws.Range("B7").Value = "Ht:`n- 800 Ht6`n- 6000 Ht4a`n- 6600 Ht7"

# Remaining column C cells (reuse existing shared strings / plain numbers).
$ws.Range("C2").Value = 3022
$ws.Range("C3").Value = "- Ehv`n- Btl"
$ws.Range("C14").Value = "Done"
$ws.Range("C15").Value = "Done"

# --- Formatting: column C mirrors column B's look -------------------------
$ws.Range("C1").Style = $ws.Range("B1").Style
$ws.Range("C2").Style = $ws.Range("B2").Style
$ws.Range("C3").Style = $ws.Range("B3").Style
$ws.Range("C4").Style = $ws.Range("B4").Style
$ws.Range("C5").Style = $ws.Range("B5").Style
$ws.Range("C6").Style = $ws.Range("B6").Style
$ws.Range("C7").Style = $ws.Range("B2").Style
$ws.Range("C8").Style = $ws.Range("B8").Style
$ws.Range("C9").Style = $ws.Range("B9").Style
$ws.Range("C10").Style = $ws.Range("B10").Style
$ws.Range("C11").Style = $ws.Range("B12").Style
$ws.Range("C12").Style = $ws.Range("B11").Style
$ws.Range("C13").Style = $ws.Range("B12").Style
$ws.Range("C14").Style = $ws.Range("B13").Style
$ws.Range("C15").Style = $ws.Range("B13").Style

$ws.Columns.Item(3).ColumnWidth = 16.17

# --- View-state refresh (selection moves to C5 / E38 as recorded) ---------
$ws.Range("C5").Select()

$qd = $wb.Worksheets.Item("QD consists")
$qd.Range("E38").Select()

$ws.Activate()
